$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header change: B1 "response" -> "verdict"
$ws.Range("B1").Value = "verdict"

# Row 3: verdict FALSE -> TRUE, confidence 0.8 -> 0.9
$ws.Cells.Item(3, 2).Value = "'TRUE"
$ws.Range("C3").Value = 0.9

# Row 6: confidence 1 -> 0.9
$ws.Range("C6").Value = 0.9

# Row 7: confidence 1 -> 0.9
$ws.Range("C7").Value = 0.9

# Row 8: confidence 0.8 -> 0.9
$ws.Range("C8").Value = 0.9

# Row 9: verdict TRUE -> FALSE, confidence 1 -> 0.8
$ws.Cells.Item(9, 2).Value = "'FALSE"
$ws.Range("C9").Value = 0.8

# Row 11: confidence 0.8 -> 0.9
$ws.Range("C11").Value = 0.9

# New row 12
$ws.Range("A12").Value = "The name of my grandfather is John"
$ws.Range("B12").Value = "INSUFFICIENT INFO"
$ws.Range("C12").Value = "'"
